$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Alcam", "Nrp1", "ECs", 2, 0.6666666666666666, 29.77446566666667, 89.323397, 0.4936227575398824, 0.4936227575398824, 3, 1, 87.94127933333334, 263.823838, 0.4109331243514438, 0.4109331243514437, 2618.404602193077, 23565.64141973769, 0.2028459420068391, 0.2028459420068391),
    @("ECs", "Alcam", "Nrp1", "FAPs", 2, 0.6666666666666666, 29.77446566666667, 89.323397, 0.4936227575398824, 0.4936227575398824, 3, 1, 52.441971, 157.325913, 0.2450515065683088, 0.2450515065683087, 1561.431665031829, 14052.88498528646, 0.1209630004115512, 0.1209630004115512),
    @("ECs", "Alcam", "Nrp1", "M2", 2, 0.6666666666666666, 29.77446566666667, 89.323397, 0.4936227575398824, 0.4936227575398824, 3, 1, 54.667459, 164.002377, 0.255450795093328, 0.255450795093328, 1627.694381079408, 14649.24942971467, 0.126096325889724, 0.126096325889724),
    @("ECs", "Alcam", "Nrp1", "sCs", 2, 0.6666666666666666, 29.77446566666667, 89.323397, 0.4936227575398824, 0.4936227575398824, 3, 1, 18.95316166666667, 56.85948500000001, 0.08856457398691947, 0.08856457398691944, 564.3202613189495, 5078.882351870546, 0.04371748923176812, 0.04371748923176811),
    @("FAPs", "Alcam", "Nrp1", "ECs", 3, 1, 0.5058753333333333, 1.517626, 0.008386769381758081, 0.008386769381758081, 3, 1, 87.94127933333334, 263.823838, 0.4109331243514438, 0.4109331243514437, 44.48732399650978, 400.3859159685881, 0.003446401345260875, 0.003446401345260874),
    @("FAPs", "Alcam", "Nrp1", "FAPs", 3, 1, 0.5058753333333333, 1.517626, 0.008386769381758081, 0.008386769381758081, 3, 1, 52.441971, 157.325913, 0.2450515065683088, 0.2450515065683087, 26.529099560282, 238.761896042538, 0.002055190472240781, 0.002055190472240781),
    @("FAPs", "Alcam", "Nrp1", "M2", 3, 1, 0.5058753333333333, 1.517626, 0.008386769381758081, 0.008386769381758081, 3, 1, 54.667459, 164.002377, 0.255450795093328, 0.255450795093328, 27.65491904411133, 248.894271397002, 0.002142406906834481, 0.00214240690683448),
    @("FAPs", "Alcam", "Nrp1", "sCs", 3, 1, 0.5058753333333333, 1.517626, 0.008386769381758081, 0.008386769381758081, 3, 1, 18.95316166666667, 56.85948500000001, 0.08856457398691947, 0.08856457398691944, 9.587936975845558, 86.29143278261002, 0.0007427706574219444, 0.0007427706574219442),
    @("M2", "Alcam", "Nrp1", "ECs", 3, 1, 29.12203633333333, 87.36610899999999, 0.482806309304492, 0.482806309304492, 3, 1, 87.94127933333334, 263.823838, 0.4109331243514438, 0.4109331243514437, 2561.029131945149, 23049.26218750634, 0.1984011051390844, 0.1984011051390844),
    @("M2", "Alcam", "Nrp1", "FAPs", 3, 1, 29.12203633333333, 87.36610899999999, 0.482806309304492, 0.482806309304492, 3, 1, 52.441971, 157.325913, 0.2450515065683088, 0.2450515065683087, 1527.216984853613, 13744.95286368252, 0.1183124134757506, 0.1183124134757506),
    @("M2", "Alcam", "Nrp1", "M2", 3, 1, 29.12203633333333, 87.36610899999999, 0.482806309304492, 0.482806309304492, 3, 1, 54.667459, 164.002377, 0.255450795093328, 0.255450795093328, 1592.02772724901, 14328.24954524109, 0.1233332555879077, 0.1233332555879077),
    @("M2", "Alcam", "Nrp1", "sCs", 3, 1, 29.12203633333333, 87.36610899999999, 0.482806309304492, 0.482806309304492, 3, 1, 18.95316166666667, 56.85948500000001, 0.08856457398691947, 0.08856457398691944, 551.9546626882072, 4967.591964193865, 0.04275953510174921, 0.04275953510174919),
    @("sCs", "Alcam", "Nrp1", "ECs", 3, 1, 0.9158823333333334, 2.747647, 0.0151841637738675, 0.0151841637738675, 3, 1, 87.94127933333334, 263.823838, 0.4109331243514438, 0.4109331243514437, 80.54386411213179, 724.8947770091861, 0.006239675860259382, 0.006239675860259382),
    @("sCs", "Alcam", "Nrp1", "FAPs", 3, 1, 0.9158823333333334, 2.747647, 0.0151841637738675, 0.0151841637738675, 3, 1, 52.441971, 157.325913, 0.2450515065683088, 0.2450515065683087, 48.030674764079, 432.2760728767111, 0.003720902208766168, 0.003720902208766168),
    @("sCs", "Alcam", "Nrp1", "M2", 3, 1, 0.9158823333333334, 2.747647, 0.0151841637738675, 0.0151841637738675, 3, 1, 54.667459, 164.002377, 0.255450795093328, 0.255450795093328, 50.06895990632434, 450.620639156919, 0.003878806708861762, 0.003878806708861761),
    @("sCs", "Alcam", "Nrp1", "sCs", 3, 1, 0.9158823333333334, 2.747647, 0.0151841637738675, 0.0151841637738675, 3, 1, 18.95316166666667, 56.85948500000001, 0.08856457398691947, 0.08856457398691944, 17.35886593131056, 156.229793381795, 0.001344778995980191, 0.00134477899598019),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}